# Add a new login entry (Admin1 / admin123) as row 2 on the LoginData
# sheet, reusing the same cell formatting (bordered style) as row 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginData")

# Copy row 1's formatting down to row 2 first ...
$ws.Range("A1:B1").Copy()
$ws.Range("A2:B2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ... then fill in the new values.
$ws.Range("A2").Value = "Admin1"
$ws.Range("B2").Value = "admin123"

$ws.Range("B6").Select() | Out-Null
